$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "condition" column (D) values used to be stored with a leading "-"
# (e.g. "-smile", "-noExistElement"). Strip the leading "-" while keeping
# the cells' existing number format / quote-prefix styling by entering the
# values with a leading apostrophe (forces text entry without altering the
# cell's style index).
$ws.Range("D2").Value = "'smile"
$ws.Range("D3").Value = "'smile"
$ws.Range("D4").Value = "'smile"
$ws.Range("D5").Value = "'smile"
$ws.Range("D6").Value = "'noExistElement"

# Move the sheet's saved selection from F7 to E8.
$ws.Range("E8").Select()
